# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New table contents for the "Estado de Cuenta" detail rows (B16:G29).
# Tipo Doc | N Doc Trabajador | Nombre Trabajador | Periodo Mora | Valor Mora | Salario Basico
$data = @(
  @("CC","1047447875","VICTOR ALFONSO VELASQUEZ IRIARTE","2010",35112,877803),
  @("CC","1108763580","MARIA CLARA BANQUET ROMERO","2109",36341,877803),
  @("CC","1108763580","MARIA CLARA BANQUET ROMERO","2110",36341,877803),
  @("CC","1108763580","MARIA CLARA BANQUET ROMERO","2111",36341,877803),
  @("CC","1108763580","MARIA CLARA BANQUET ROMERO","2112",36341,877803),
  @("CC","1108763580","MARIA CLARA BANQUET ROMERO","2201",36341,877803),
  @("CC","1143348923","PABLO EMILIO LONDOÑO RIOS","2303",40000,1000000),
  @("CC","1127587489","KAREN CECILIA BARRAGAN MUNZON","2303",46400,1300000),
  @("CC","1143391729","EDWIN MANUEL MARTINEZ LOZANO","2303",46400,1000000),
  @("CC","1143348923","PABLO EMILIO LONDOÑO RIOS","2304",40000,1000000),
  @("CC","1127587489","KAREN CECILIA BARRAGAN MUNZON","2304",46400,1300000),
  @("CC","1143391729","EDWIN MANUEL MARTINEZ LOZANO","2304",46400,1000000),
  @("CC","1143348923","PABLO EMILIO LONDOÑO RIOS","2309",24000,1000000),
  @("CC","1127587489","KAREN CECILIA BARRAGAN MUNZON","2309",31200,1300000)
)

$r = 16
foreach ($row in $data) {
  $ws.Cells.Item($r, 2).Value = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
  $ws.Cells.Item($r, 6).Value = $row[4]
  $ws.Cells.Item($r, 7).Value = $row[5]
  $r++
}
